$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarterly indexing bug-fix: the per-quarter error-statistics rows had been
# off by one. Shift the existing B:G data for rows 2-10 down into rows 3-11
# (the row's quarter index in column A stays put), and populate row 2 with
# the newly computed statistics for the quarter that was previously missing.

# Capture current B:G values for rows 2 through 10 before overwriting anything.
$data = @()
for ($r = 2; $r -le 10; $r++) {
    $row = @()
    for ($c = 2; $c -le 7; $c++) {
        $row += $ws.Cells.Item($r, $c).Value2
    }
    $data += , $row
}

# Shift rows 2-10 down to rows 3-11.
for ($i = 0; $i -lt $data.Count; $i++) {
    $destRow = 3 + $i
    $row = $data[$i]
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $row[$c - 2]
    }
}

# Write the newly computed statistics into row 2.
$ws.Cells.Item(2, 2).Value2 = -0.02314597604078636
$ws.Cells.Item(2, 3).Value2 = 0.3579920056255013
$ws.Cells.Item(2, 4).Value2 = 0.1782699060034266
$ws.Cells.Item(2, 5).Value2 = 0.4222202103209018
$ws.Cells.Item(2, 6).Value2 = 0.4363822494547141
$ws.Cells.Item(2, 7).Value2 = 15
